$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Clear all existing contents on both sheets (rows 1-13)
$ws1.Range("A1:A13").ClearContents()
$ws2.Range("A1:A13").ClearContents()

# Sheet 1 ("Worksheet"): new shortened data - a single label followed by numbers
$ws1.Range("A1").Value = "ABM"
$ws1.Range("A2").Value = 2
$ws1.Range("A3").Value = 3
$ws1.Range("A4").Value = 4
$ws1.Range("A5").Value = 5

# Sheet 2 ("Worksheet 1"): keep cells A1:A5 present but empty
$ws2.Range("A1:A5").Borders.LineStyle = 0
